# msz - table part 1
#
# Adds a new "Select <Tier>" step before each existing "Choose <Tier>" step
# in the product-insurance test table: 4 new rows are inserted at row 15
# (pushing the existing "Choose Silver/Gold/Platinum/Ultimate" and
# "Send Quote" rows down by 4), and the new rows are populated with
# "Select Silver" / "Select Gold" / "Select Platinum" / "Select Ultimate".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the existing "Choose Silver" row (row 15),
# shifting everything below (including the floating picture) down.
$ws.Rows("15:18").Insert()

# Row 15 - Select Silver
$ws.Range("A15").Value = "Select Silver"
$ws.Range("B15").Value = "<SET>"
$ws.Range("F15").Value = "Select Silver"
$ws.Range("H15").Value = "<NOP>"

# Row 16 - Select Gold
$ws.Range("A16").Value = "Select Gold"
$ws.Range("B16").Value = "<SET>"
$ws.Range("F16").Value = "Select Gold"
$ws.Range("H16").Value = "<NOP>"

# Row 17 - Select Platinum
$ws.Range("A17").Value = "Select Platinum"
$ws.Range("B17").Value = "<SET>"
$ws.Range("F17").Value = "Select Platinum"
$ws.Range("H17").Value = "<NOP>"

# Row 18 - Select Ultimate
$ws.Range("A18").Value = "Select Ultimate"
$ws.Range("B18").Value = "<SET>"
$ws.Range("F18").Value = "Select Ultimate"
$ws.Range("H18").Value = "<NOP>"

# The picture anchored below the table (row 15-18 inserted above it) needs
# to move down along with the rows it used to sit below. Re-anchor it to
# its new absolute position on the sheet.
$shp = $ws.Shapes.Item(1)
$shp.Top = 366.0000787401575

# Restore the active-cell selection to match where editing left off.
$ws.Range("F28").Select() | Out-Null
